# Rebuilt the df module.
# Turns the generic Sheet1 (a/b/c demo data) into a two-sheet "settings" /
# "tasks" workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "settings" (reuse the original Sheet1 so sheetId/rId stay at 1)
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item(1)
$settings.Name = "settings"

# A1 already carries the bordered/bold "label" style (s=1) in the source
# workbook - grab a copy of that formatting before we overwrite the cell,
# so we can stamp it onto the new label cells further down.
$settings.Range("A1").Copy()

$settings.Range("B1").Value = "value"
$settings.Range("C1").Value = "help"

$settings.Range("A2").Value = "height"
$settings.Range("B2").Value = 100
$settings.Range("C2").Value = "The height of the object"

$settings.Range("A3").Value = "width"
$settings.Range("B3").Value = 200
$settings.Range("C3").Value = "The width of the object"

$settings.Range("A4").Value = "rows"
$settings.Range("B4").Value = 5
$settings.Range("C4").Value = "The number of rows in the object"

# Stamp the label style onto the row-header cells in column A (B1/C1 keep
# their original style automatically since we only changed their value).
$settings.Range("A2:A4").PasteSpecial(-4122)

# A1 itself is no longer used - drop its value/format.
$settings.Range("A1").ClearContents()
$settings.Range("A1").ClearFormats()

$settings.Range("A1").Select()

# ---------------------------------------------------------------------
# Sheet "tasks" - duplicate "settings" (so it inherits the same
# sheetFormatPr / pageMargins / namespaces) and wipe the copy clean.
# ---------------------------------------------------------------------
$settings.Copy([System.Reflection.Missing]::Value, $settings)
$tasks = $wb.Worksheets.Item(2)
$tasks.Name = "tasks"
$tasks.Cells.Clear()

# Re-use the same label style for the header row / task-index column.
$settings.Range("B1").Copy()

$tasks.Range("B1").Value = "task"
$tasks.Range("C1").Value = "start"
$tasks.Range("D1").Value = "finish"

$tasks.Range("A2").Value = 0
$tasks.Range("B2").Value = "task1"
$tasks.Range("C2").Value = 45366
$tasks.Range("D2").Value = 45371

$tasks.Range("A3").Value = 1
$tasks.Range("B3").Value = "task2"
$tasks.Range("C3").Value = 45376
$tasks.Range("D3").Value = 45381

$tasks.Range("A4").Value = 2
$tasks.Range("B4").Value = "task3"
$tasks.Range("C4").Value = 45392
$tasks.Range("D4").Value = 45402

$tasks.Range("B1:D1").PasteSpecial(-4122)
$tasks.Range("A2:A4").PasteSpecial(-4122)

# Date columns get a plain "d-mmm-yy" (built-in numFmtId 15) number format.
$tasks.Range("C2:D4").NumberFormat = "d-mmm-yy"

$tasks.Columns("C:D").AutoFit()

$tasks.Range("C4").Select()

Write-Host "done"
